$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Marc"

# "2026-02-11" looks like a date, so Excel would normally auto-convert it
# to a date serial number. Force the cell to Text first so the literal
# string is stored verbatim (matching the source row's inline string),
# then reset the style back to Normal so no extra formatting sticks.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2026-02-11"
$ws.Range("B3").Style = "Normal"

$ws.Range("C3").Value = "16:59:47"
